# Commit: "update MLP model training (add parameter about MLP structure)"
#
# This adds a new "mlp_structure" parameter block to the
# config_file_structure sheet, right after the closing "}," of the
# "model_type" section and before the "dnn_training_parameter" block.
# Three new rows are inserted (old row 28 -> new row 31, everything
# below shifts down by 3 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config_file_structure")
$ws.Activate()

# Insert 3 blank rows before the old row 28 ("dnn_training_parameter": {)
$ws.Range("A28:A30").EntireRow.Insert()

# Make sure the new rows have the same row height/format as their neighbours
$ws.Rows.Item("28:30").RowHeight = 18.75

# Row 28: "mlp_structure": {      // structure parameters (show on MLP model only)
$ws.Cells.Item(28, 4).Value = """mlp_structure"": {"
$ws.Cells.Item(28, 20).Value = "// structure parameters (show on MLP model only)"

# Row 29: "num_of_hidden_nodes": <parameter block>,      // number of hidden nodes (configurable)
# (column D has no entry on this row, clear the auto-populated blank cell so
# it doesn't linger as an empty styled cell)
$ws.Cells.Item(29, 4).Clear()
$ws.Cells.Item(29, 5).Value = """num_of_hidden_nodes"": <parameter block>,"
$ws.Cells.Item(29, 20).Value = "// number of hidden nodes (configurable)"

# Row 30: },
$ws.Cells.Item(30, 4).Value = "},"

# Restore the selection recorded in the workbook (P30) on the active sheet
$ws.Range("P30").Select()
